## Data Overview.xlsx — "load and combine TCR data for Dataset 1 (Adrienne)
## with scRepertoire"
##
## The "Part 1 (22 Samples)" sheet listed only one "+CPI no colitis" sample
## (NC1) repeated down column B. This adds the remaining Dataset-1 samples
## (NC2..NC6) to rows 3-7 of that column, and leaves the user focused on
## that sheet (it becomes the active tab / selected sheet with C2 selected)
## instead of the previously-active "Part 5" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Part 1 (22 Samples)")

# Column B currently holds shared-string index 0 ("+CPI no colitis NC1")
# for every data row. Rows 3-7 should instead hold NC2..NC6.
$ws1.Range("B3").Value = '"+CPI no colitis NC2"'
$ws1.Range("B4").Value = '"+CPI no colitis NC3"'
$ws1.Range("B5").Value = '"+CPI no colitis NC4"'
$ws1.Range("B6").Value = '"+CPI no colitis NC5"'
$ws1.Range("B7").Value = '"+CPI no colitis NC6"'

# Make "Part 1 (22 Samples)" the active sheet/tab (it was "Part 5" before),
# with C2 as the selected cell.
[void]$ws1.Activate()
$ws1.Range("C2").Select() | Out-Null
